$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation is inserted above the current row 188, pushing
# every following row (188-260) down by one (to 189-261). The new row 188
# re-uses the same descriptive fields (market/region/product/quality/unit/
# origin/box-size) that the old row 188 had, but carries its own date and
# price/volume figures.

$ws.Rows("188:188").Insert()

$newRow = 188

$ws.Cells.Item($newRow, 1).Value = 7
$ws.Cells.Item($newRow, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($newRow, 3).Value = "Ñuble"
$ws.Cells.Item($newRow, 4).Value = 44875
$ws.Cells.Item($newRow, 5).Value = 16
$ws.Cells.Item($newRow, 6).Value = "Fruta"
$ws.Cells.Item($newRow, 7).Value = 100108
$ws.Cells.Item($newRow, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($newRow, 9).Value = 100108005
$ws.Cells.Item($newRow, 10).Value = "Piña"
$ws.Cells.Item($newRow, 11).Value = "Caramelo"
$ws.Cells.Item($newRow, 12).Value = "Segunda"
$ws.Cells.Item($newRow, 13).Value = 30
$ws.Cells.Item($newRow, 14).Value = 29000
$ws.Cells.Item($newRow, 15).Value = 30000
$ws.Cells.Item($newRow, 16).Value = 29500
$ws.Cells.Item($newRow, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item($newRow, 18).Value = "Ecuador"
$ws.Cells.Item($newRow, 19).Value = 2107
$ws.Cells.Item($newRow, 20).Value = 14

# Keep the date column's display format consistent with the rest of the
# column (yyyy-mm-dd hh:mm:ss), matching what Insert() already copied down
# from the row above, but set explicitly to be safe.
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
